$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dist_coûts")

# Update cost values (maj plans + cout alpha III)
$ws.Range("C2").Value = 1.8
$ws.Range("M2").Value = 30
$ws.Range("C3").Value = 0.67
$ws.Range("M3").Value = 30
$ws.Range("C4").Value = 0.88
$ws.Range("C5").Value = 1.8
$ws.Range("C6").Value = 0.16
$ws.Range("C7").Value = 0.23
$ws.Range("C8").Value = 0.62

# Move active selection to M4
$ws.Range("M4").Select()
